$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per the PCA re-run (7 PCs)
$ws.Range("B2").Value = 154
$ws.Range("B3").Value = 131
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 30

# Remove the now-unused last row (previously row 5: A5=2, B5=50)
$ws.Rows.Item(5).Delete()
